# Apply the graphical updates:
#  - Fix park-name labels (shared strings): a typo "ParkEiffel" becomes the
#    correctly spelled "ParkSaechs_Schw", and a missing park "ParkEifel" is
#    inserted, which shifts the intervening labels (ParkHainich,
#    ParkHunsrueck) down by one row.
#  - The numeric model-summary rows (rows 5,6,7,10 = Value/Std.Error/p for
#    columns B,C,D) follow the same row rotation as the labels, on every
#    one of the 8 worksheets.

$wb = $excel.ActiveWorkbook

# New label text for rows 5..10 (same layout on every sheet).
$labels = @{
    5  = "ParkEifel"
    6  = "ParkHainich"
    7  = "ParkHunsrueck"
    8  = "ParkJasmund"
    9  = "ParkKellerwald"
    10 = "ParkSaechs_Schw"
}

# New B/C/D values for rows 5,6,7,10 per sheet (1-indexed sheet order),
# obtained by rotating the old row5/row6/row7/row10 values: the value that
# used to sit in row5 now sits in row6, row6 -> row7, row7 -> row10, and
# row10 -> row5. Rows 8 and 9 are untouched.
$data = @{
    1 = @{ 5 = @(0.1788, 0.368, 0.6271); 6 = @(0.2012, 0.3619, 0.5782); 7 = @(-4.5034, 76547.695, 1); 10 = @(-4.4454, 67350.3807, 0.9999) }
    2 = @{ 5 = @(-0.304, 0.1885, 0.1069); 6 = @(0.2428, 0.1882, 0.1971); 7 = @(-0.4781, 0.2056, 0.0201); 10 = @(0.4481, 0.1745, 0.0102) }
    3 = @{ 5 = @(0.5337, 0.1185, 0); 6 = @(-0.2097, 0.1206, 0.0822); 7 = @(0.3294, 0.1316, 0.0123); 10 = @(0.1474, 0.1129, 0.1916) }
    4 = @{ 5 = @(-0.0065, 0.0706, 0.9264); 6 = @(-0.1128, 0.0824, 0.171); 7 = @(0.0491, 0.0715, 0.4922); 10 = @(-0.1129, 0.071, 0.1119) }
    5 = @{ 5 = @(0.2896, 0.1193, 0.0153); 6 = @(0.0157, 0.1312, 0.9046); 7 = @(0.1359, 0.1298, 0.2948); 10 = @(-0.1722, 0.1279, 0.1783) }
    6 = @{ 5 = @(0.3727, 0.1055, 0.0004); 6 = @(-0.1849, 0.1226, 0.1315); 7 = @(0.4938, 0.1076, 0); 10 = @(0.1287, 0.1034, 0.2132) }
    7 = @{ 5 = @(-0.1254, 0.259, 0.6284); 6 = @(0.1766, 0.2343, 0.451); 7 = @(-0.3861, 0.2777, 0.1644); 10 = @(-0.6188, 0.266, 0.02) }
    8 = @{ 5 = @(0.0157, 0.1299, 0.9036); 6 = @(-0.0679, 0.1331, 0.6098); 7 = @(0.0203, 0.1366, 0.882); 10 = @(-0.3214, 0.1248, 0.01) }
}

$sheetCount = $wb.Worksheets.Count
for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Update the park labels in column A.
    foreach ($row in $labels.Keys) {
        $ws.Cells.Item($row, 1).Value = $labels[$row]
    }

    # Update the rotated numeric values in columns B, C, D.
    $sheetData = $data[$i]
    foreach ($row in $sheetData.Keys) {
        $vals = $sheetData[$row]
        $ws.Cells.Item($row, 2).Value = $vals[0]
        $ws.Cells.Item($row, 3).Value = $vals[1]
        $ws.Cells.Item($row, 4).Value = $vals[2]
    }
}
